# Add new day of data ("20250906") for both cohorts:
#  - RZJ_3M: day10 (5 sessions) inserted right after the existing RZJ_3M block (after row 36)
#  - RZK_15M: day6 (5 sessions) appended after the existing RZK_15M block (at the end of the sheet)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# 1) Insert 5 new rows right after the last RZJ_3M row (row 36), pushing the
#    RZK_15M block down by 5 rows (was 37-61, becomes 42-66).
# ---------------------------------------------------------------------------
for ($i = 0; $i -lt 5; $i++) {
    $ws.Rows.Item(37).Insert()
}

# Copy formatting (fonts/styles) from the last RZJ_3M row (row 36) into the
# freshly inserted rows 37-41 so they pick up the same look (blue "Mice" font
# for RZJ_3M, green/red status columns, etc.)
$ws.Rows.Item(36).Copy()
$ws.Range("A37:N41").PasteSpecial($xlPasteFormats)

# Fill in the new RZJ_3M day10 rows (sessions 1-5)
$rzjDay10 = 37
for ($s = 1; $s -le 5; $s++) {
    $r = $rzjDay10 + ($s - 1)
    $ws.Cells.Item($r, 1).Value = "202508_Preliminary"
    $ws.Cells.Item($r, 2).Value = "RZJ_3M"
    $ws.Cells.Item($r, 3).Value = "Default"
    $ws.Cells.Item($r, 4).Value = "20250906"
    $ws.Cells.Item($r, 5).Value = "20250906_RZJ_3M_Thermal_day10_" + $s + "_PelDisLeft50passive"
    $ws.Cells.Item($r, 6).Value = "✅"
    $ws.Cells.Item($r, 7).Value = "✅"
    $ws.Cells.Item($r, 8).Value = "✅"
    $ws.Cells.Item($r, 9).Value = "❌"
    $ws.Cells.Item($r, 10).Value = "❌"
    $ws.Cells.Item($r, 11).Value = "❌"
    $ws.Cells.Item($r, 12).Value = "❌"
    $ws.Cells.Item($r, 13).Value = "❌"
    $ws.Cells.Item($r, 14).Value = "✅"
}

# ---------------------------------------------------------------------------
# 2) Append 5 new rows after the (now shifted) last RZK_15M row, which is
#    row 66 (was row 61 before the insert above), for RZK_15M day6.
# ---------------------------------------------------------------------------
$lastRzkRow = 66

$ws.Rows.Item($lastRzkRow).Copy()
$ws.Range("A67:N71").PasteSpecial($xlPasteFormats)

$rzkDay6 = $lastRzkRow + 1
for ($s = 1; $s -le 5; $s++) {
    $r = $rzkDay6 + ($s - 1)
    $ws.Cells.Item($r, 1).Value = "202508_Preliminary"
    $ws.Cells.Item($r, 2).Value = "RZK_15M"
    $ws.Cells.Item($r, 3).Value = "default"
    $ws.Cells.Item($r, 4).Value = "20250906"
    $ws.Cells.Item($r, 5).Value = "20250906_RZK_15M_Thermal_day6_" + $s + "_PelSAT50passive"
    $ws.Cells.Item($r, 6).Value = "✅"
    $ws.Cells.Item($r, 7).Value = "✅"
    $ws.Cells.Item($r, 8).Value = "✅"
    $ws.Cells.Item($r, 9).Value = "❌"
    $ws.Cells.Item($r, 10).Value = "❌"
    $ws.Cells.Item($r, 11).Value = "❌"
    $ws.Cells.Item($r, 12).Value = "❌"
    $ws.Cells.Item($r, 13).Value = "❌"
    $ws.Cells.Item($r, 14).Value = "✅"
}

Write-Output "Added RZJ_3M day10 (rows 37-41) and RZK_15M day6 (rows 67-71)."
